$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing row 100: correct the date/time value in column A ---
$ws.Range("A100").Value = 45467.2916666667

# --- Append new row 101 with data from the R script results ---

# Column A: date (apply same style/number format as A100 by copying it first)
$ws.Range("A100").Copy()
$ws.Range("A101").PasteSpecial(-4122) | Out-Null
$ws.Range("A101").Value = 45468.2916666667

# Column B: volume (plain number)
$ws.Range("B101").Value = 13221

# Columns C-F: high, low, open, close (plain numbers)
$ws.Range("C101").Value = 0.735000014305115
$ws.Range("D101").Value = 0.704999983310699
$ws.Range("E101").Value = 0.709999978542328
$ws.Range("F101").Value = 0.704999983310699

# Column G: adj_close stored as text (matches source data / shared string "0.704999983310699")
# Use a Text number format temporarily so the numeric-looking string isn't coerced to a
# number, then restore the General formatting (copied from an existing plain-text cell)
# so no stray style is left behind on the cell.
$ws.Range("G101").NumberFormat = "@"
$ws.Range("G101").Value = "0.704999983310699"
$ws.Range("G100").Copy()
$ws.Range("G101").PasteSpecial(-4122) | Out-Null

# Column H: ticker (plain text, matches existing shared string "BWZ.MI")
$ws.Range("H101").Value = "BWZ.MI"
